$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) and "Volume(1h)" column (E) hold plain text in the
# source workbook (t="inlineStr"), even though many Price values look like
# numbers (thousand-grouping dots make them invalid numerics anyway, e.g.
# "56.579.02"). Whenever a new value would otherwise be auto-converted to a
# number by Excel, prefix it with a leading apostrophe so it is stored as
# literal text (quote-prefixed), matching the original cell content exactly.

$ws.Range("D2").Value = "56.579.02"
$ws.Range("E2").Value = "  -2.65%  "

$ws.Range("D3").Value = "2.975.99"
$ws.Range("E3").Value = "  -4.93%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'495.81"
$ws.Range("E5").Value = "  -5.45%  "

$ws.Range("D6").Value = "'134.19"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "2.973.94"
$ws.Range("E8").Value = "  -4.96%  "

$ws.Range("D9").Value = "'0.425"
$ws.Range("E9").Value = "  -4.30%  "

$ws.Range("D10").Value = "'7.21"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("E11").Value = "  -2.91%  "

$ws.Range("D12").Value = "'0.349"
$ws.Range("E12").Value = "  -6.33%  "

$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("D14").Value = "3.491.84"
$ws.Range("E14").Value = "  -4.76%  "

$ws.Range("D15").Value = "'25.16"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("D16").Value = "56.550.64"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").Value = "'0.0000147"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").Value = "2.982.20"
$ws.Range("E18").Value = "  -4.76%  "

$ws.Range("D19").Value = "'5.69"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'12.26"
$ws.Range("E20").Value = "  -5.22%  "

$ws.Range("D21").Value = "'7.73"
$ws.Range("E21").Value = "  -1.70%  "

$ws.Range("D22").Value = "'324.86"
$ws.Range("E22").Value = "  -4.98%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "'0.467"
$ws.Range("E24").Value = "  -8.19%  "

$ws.Range("D25").Value = "'61.76"
$ws.Range("E25").Value = "  -8.07%  "

$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").Value = "'0.162"
$ws.Range("E27").Value = "  -4.37%  "

$ws.Range("D28").Value = "0.0₃0889"
$ws.Range("E28").Value = "  -3.75%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").Value = "'6.36"
$ws.Range("E30").Value = "  -5.46%  "

$ws.Range("D31").Value = "'6.76"
$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  -8.08%  "

$ws.Range("D33").Value = "'20.25"
$ws.Range("E33").Value = "  -5.26%  "

$ws.Range("E34").Value = "  -7.29%  "

$ws.Range("D35").Value = "'153.04"
$ws.Range("E35").Value = "  -4.11%  "

$ws.Range("D36").Value = "'4.41"
$ws.Range("E36").Value = "  -8.05%  "

$ws.Range("E37").Value = "  -7.40%  "

$ws.Range("D38").Value = "'5.56"
$ws.Range("E38").Value = "  -10.11%  "

$ws.Range("D39").Value = "'0.0670"
$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").Value = "'22.97"
$ws.Range("E40").Value = "  -2.35%  "

$ws.Range("D41").Value = "3.010.63"
$ws.Range("E41").Value = "  -4.75%  "

$ws.Range("D42").Value = "'36.36"
$ws.Range("E42").Value = "  -9.98%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'0.638"
$ws.Range("E44").Value = "  -7.05%  "

$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -7.28%  "

$ws.Range("D46").Value = "2.220.36"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").Value = "  -3.50%  "

$ws.Range("D48").Value = "'3.53"
$ws.Range("E48").Value = "  -9.82%  "

$ws.Range("D49").Value = "'1.93"
$ws.Range("E49").Value = "  +7.09%  "

$ws.Range("E50").Value = "  +1.06%  "

$ws.Range("D51").Value = "'5.70"
$ws.Range("E51").Value = "  -7.07%  "

